{"js": "// Remove the \"personal message functionality\" bullet (messaging between\n// members is no longer offered) and keep the trailing _GoBack bookmark on\n// the paragraph that is left behind, so the final blank paragraph at the\n// end of the document no longer carries it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the bullet describing the personal-messaging feature.\nconst target = paragraphs.items.find(p =>\n    p.text.indexOf(\"personal message functionality\") !== -1\n);\n\nif (target) {\n    target.delete();\n    await context.sync();\n}\n\n// Move the \"_GoBack\" bookmark off the final (now trailing empty) paragraph\n// and onto the blank paragraph that used to sit right after the deleted\n// bullet (immediately before \"Non-Functional Requirements\").\nconst bookmarkRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nbookmarkRange.load(\"isNullObject\");\nawait context.sync();\n\nif (!bookmarkRange.isNullObject) {\n    context.document.deleteBookmark(\"_GoBack\");\n\n    const refreshedParagraphs = body.paragraphs;\n    refreshedParagraphs.load(\"items/text,items/style,items/isLastParagraph\");\n    await context.sync();\n\n    // Prefer the blank \"List Paragraph\"-styled paragraph (the spacer that sat\n    // right after the deleted bullet); fall back to the first blank,\n    // non-final paragraph if styles ever differ.\n    let blankBeforeNonFunctional = refreshedParagraphs.items.find(p =>\n        p.text.trim() === \"\" && p.style === \"List Paragraph\" && !p.isLastParagraph\n    );\n    if (!blankBeforeNonFunctional) {\n        blankBeforeNonFunctional = refreshedParagraphs.items.find(p =>\n            p.text.trim() === \"\" && !p.isLastParagraph\n        );\n    }\n\n    if (blankBeforeNonFunctional) {\n        blankBeforeNonFunctional.getRange().insertBookmark(\"_GoBack\");\n    }\n\n    await context.sync();\n}\n", "ps1": "# Remove the \"personal message functionality\" bullet (messaging between\n# members is no longer offered) and re-home the trailing \"_GoBack\" bookmark\n# onto the blank paragraph left behind, so the final paragraph of the\n# document no longer carries it.\n\n$d = $word.ActiveDocument\n\n# Locate and remove the paragraph describing the personal-messaging feature.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*personal message functionality*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# Drop the existing \"_GoBack\" bookmark (it currently sits on the final,\n# trailing empty paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Re-create the \"_GoBack\" bookmark on the blank paragraph that immediately\n# precedes \"Non-Functional Requirements\" (the paragraph left behind once the\n# messaging bullet was removed).\n$count2 = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count2; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t.Trim() -eq \"\" -and $i -lt $count2) {\n        $nextp = $d.Paragraphs.Item($i + 1)\n        $nextt = $nextp.Range.Text\n        if ($nextt -like \"*Non-Functional Requirements*\") {\n            $d.Bookmarks.Add(\"_GoBack\", $p.Range)\n            break\n        }\n    }\n}\n"}
